$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for columns A..T
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Full replacement data for rows 2..17 (Sending cluster x Target cluster = ECs, FAPs, M2, sCs)
# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#          Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
#          Ligand total expression value, Ligand derived specificity of average expression value,
#          Ligand derived specificity of total expression value, Receptor-expressing cells,
#          Receptor detection rate, Receptor average expression value, Receptor total expression value,
#          Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
#          Edge average expression weight, Edge total expression weight,
#          Edge average expression derived specificity, Edge total expression derived specificity
$data = @(
    @("ECs", "Ltbp1", "Itgb5", "ECs", 3, 1, 6.334389333333334, 19.003168, 0.1087186954512644, 0.1087186954512644, 3, 1, 9.460162333333333, 28.380487, 0.08197024919772995, 0.08197024919772995, 59.92435137586845, 539.319162382816, 0.008911698558592255, 0.008911698558592253),
    @("ECs", "Ltbp1", "Itgb5", "FAPs", 3, 1, 6.334389333333334, 19.003168, 0.1087186954512644, 0.1087186954512644, 3, 1, 46.18256633333333, 138.547699, 0.4001618933742075, 0.4001618933742075, 292.5383555678258, 2632.845200110432, 0.04350507901695181, 0.0435050790169518),
    @("ECs", "Ltbp1", "Itgb5", "M2", 3, 1, 6.334389333333334, 19.003168, 0.1087186954512644, 0.1087186954512644, 3, 1, 51.40166966666666, 154.205009, 0.4453842886934318, 0.4453842886934319, 325.5981880520569, 2930.383692468512, 0.04842159884123925, 0.04842159884123925),
    @("ECs", "Ltbp1", "Itgb5", "sCs", 3, 1, 6.334389333333334, 19.003168, 0.1087186954512644, 0.1087186954512644, 3, 1, 8.365307333333334, 25.095922, 0.07248356873463072, 0.07248356873463073, 52.98911354232179, 476.9020218808961, 0.007880319034481107, 0.007880319034481109),
    @("FAPs", "Ltbp1", "Itgb5", "ECs", 3, 1, 25.48320566666667, 76.449617, 0.4373745802799199, 0.4373745802799199, 3, 1, 9.460162333333333, 28.380487, 0.08197024919772995, 0.08197024919772995, 241.0752623803865, 2169.677361423479, 0.03585170333829758, 0.03585170333829758),
    @("FAPs", "Ltbp1", "Itgb5", "FAPs", 3, 1, 25.48320566666667, 76.449617, 0.4373745802799199, 0.4373745802799199, 3, 1, 46.18256633333333, 138.547699, 0.4001618933742075, 0.4001618933742075, 1176.879836086809, 10591.91852478128, 0.1750206401585621, 0.1750206401585621),
    @("FAPs", "Ltbp1", "Itgb5", "M2", 3, 1, 25.48320566666667, 76.449617, 0.4373745802799199, 0.4373745802799199, 3, 1, 51.40166966666666, 154.205009, 0.4453842886934318, 0.4453842886934319, 1309.879319725728, 11788.91387753155, 0.1947997663305604, 0.1947997663305605),
    @("FAPs", "Ltbp1", "Itgb5", "sCs", 3, 1, 25.48320566666667, 76.449617, 0.4373745802799199, 0.4373745802799199, 3, 1, 8.365307333333334, 25.095922, 0.07248356873463072, 0.07248356873463073, 213.1748472402082, 1918.573625161874, 0.03170247045249984, 0.03170247045249984),
    @("M2", "Ltbp1", "Itgb5", "ECs", 2, 0.6666666666666666, 0.037545, 0.112635, 0.0006443941484995116, 0.0006443941484995116, 3, 1, 9.460162333333333, 28.380487, 0.08197024919772995, 0.08197024919772995, 0.355181794805, 3.196636153245, 0.00005282114893406396, 0.00005282114893406396),
    @("M2", "Ltbp1", "Itgb5", "FAPs", 2, 0.6666666666666666, 0.037545, 0.112635, 0.0006443941484995116, 0.0006443941484995116, 3, 1, 46.18256633333333, 138.547699, 0.4001618933742075, 0.4001618933742075, 1.733924452985, 15.605320076865, 0.0002578619825428248, 0.0002578619825428248),
    @("M2", "Ltbp1", "Itgb5", "M2", 2, 0.6666666666666666, 0.037545, 0.112635, 0.0006443941484995116, 0.0006443941484995116, 3, 1, 51.40166966666666, 154.205009, 0.4453842886934318, 0.4453842886934319, 1.929875687635, 17.368881188715, 0.0002870030294676647, 0.0002870030294676647),
    @("M2", "Ltbp1", "Itgb5", "sCs", 2, 0.6666666666666666, 0.037545, 0.112635, 0.0006443941484995116, 0.0006443941484995116, 3, 1, 8.365307333333334, 25.095922, 0.07248356873463072, 0.07248356873463073, 0.31407546383, 2.826679174470001, 0.00004670798755495818, 0.00004670798755495819),
    @("sCs", "Ltbp1", "Itgb5", "ECs", 3, 1, 26.40889, 79.22667000000001, 0.4532623301203161, 0.4532623301203161, 3, 1, 9.460162333333333, 28.380487, 0.08197024919772995, 0.08197024919772995, 249.8323864431433, 2248.49147798829, 0.03715402615190605, 0.03715402615190605),
    @("sCs", "Ltbp1", "Itgb5", "FAPs", 3, 1, 26.40889, 79.22667000000001, 0.4532623301203161, 0.4532623301203161, 3, 1, 46.18256633333333, 138.547699, 0.4001618933742075, 0.4001618933742075, 1219.630314214704, 10976.67282793233, 0.1813783122161507, 0.1813783122161507),
    @("sCs", "Ltbp1", "Itgb5", "M2", 3, 1, 26.40889, 79.22667000000001, 0.4532623301203161, 0.4532623301203161, 3, 1, 51.40166966666666, 154.205009, 0.4453842886934318, 0.4453842886934319, 1357.461040043337, 12217.14936039003, 0.2018759204921645, 0.2018759204921645),
    @("sCs", "Ltbp1", "Itgb5", "sCs", 3, 1, 26.40889, 79.22667000000001, 0.4532623301203161, 0.4532623301203161, 3, 1, 8.365307333333334, 25.095922, 0.07248356873463072, 0.07248356873463073, 220.9184811821934, 1988.266330639741, 0.03285407126009481, 0.03285407126009481)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $rowNum = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value = $rowVals[$c]
    }
}
